$wb = $excel.ActiveWorkbook

# Compute the BGR-packed color value Excel/COM expects for the standard
# "HyperLink" font color FF6495ED (stored as RGB in the OOXML font def).
$hyperlinkColor = 0x64 + (0x95 * 256) + (0xED * 65536)

function Set-HandbackRow8($SheetName, $XlfFileName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # I8: Latest Target File -> link to the source .md file (same md file
    # the A8 hyperlink already points at).
    $ws.Range("I8").Value = "fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.md"
    $ws.Hyperlinks.Add(
        $ws.Range("I8"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/928ff6336e08c8e28e68011065e0363ce95ea1e8/e2e/fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.md",
        [System.Reflection.Missing]::Value,
        [System.Reflection.Missing]::Value,
        "fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.md"
    ) | Out-Null
    $ws.Range("I8").Font.Name = "Calibri"
    $ws.Range("I8").Font.Underline = 2
    $ws.Range("I8").Font.Color = $hyperlinkColor

    # J8: Latest Handback File
    $ws.Range("J8").Value = $XlfFileName

    # K8: Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # P8: Error Detail
    $ws.Range("P8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57a08066b4ec53c068bcf29ed51b7f90c5f8255e/e2e/fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/928ff6336e08c8e28e68011065e0363ce95ea1e8/e2e/fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.md."

    # Column P (16) grows to fit the new, much longer Error Detail text.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Set-HandbackRow8 "zh-cn" "fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.8a0a45205117979195587e8eec092d47430b6680.zh-cn.xlf" "2016-08-23 06:41:14"
Set-HandbackRow8 "de-de" "fa6cf42d-5ccd-4002-9a0a-6f5419a702b3.8a0a45205117979195587e8eec092d47430b6680.de-de.xlf" "2016-08-23 06:41:21"
